# Implement relative rubric type.
# Adds a new "Relative Samples" worksheet between "Formula Samples" and
# "Test Case Samples", populates it with sample data, and makes it the
# active/selected sheet (mirroring the author's edit).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after "Formula Samples" (i.e. right before
# "Test Case Samples"), matching the target sheet order:
#   Constant Samples, Formula Samples, Relative Samples, Test Case Samples
$formulaSheet = $wb.Worksheets.Item("Formula Samples")
$relativeSheet = $wb.Worksheets.Add($null, $formulaSheet)
$relativeSheet.Name = "Relative Samples"

# Populate the sample data (no header row - data starts on row 2).
$relativeSheet.Range("A2").Value = 3003
$relativeSheet.Range("B2").Value = 1001
$relativeSheet.Range("C2").Value = 2002

$relativeSheet.Range("A3").Value = 4005
$relativeSheet.Range("B3").Value = 2002
$relativeSheet.Range("C3").Value = 2002

# Match the author's saved selection/active cell on the new sheet, and make
# it the active (tab-selected) sheet.
$relativeSheet.Activate()
$relativeSheet.Range("L7").Select()
